$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.73255033333334
$ws.Range("H2").Value = 191.197651
$ws.Range("I2").Value = 0.09718402276460011
$ws.Range("J2").Value = 0.1059076069828809
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 42.05160772145455
$ws.Range("R2").Value = 378.464469493091
$ws.Range("S2").Value = 0.0007487771762408797
$ws.Range("T2").Value = 0.0008985919980662878
$ws.Range("G3").Value = 63.73255033333334
$ws.Range("H3").Value = 191.197651
$ws.Range("I3").Value = 0.09718402276460011
$ws.Range("J3").Value = 0.1059076069828809
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("Q3").Value = 3891.51967807302
$ws.Range("R3").Value = 35023.67710265718
$ws.Range("S3").Value = 0.06929297769385123
$ws.Range("T3").Value = 0.08315706895671944
$ws.Range("G4").Value = 63.73255033333334
$ws.Range("H4").Value = 191.197651
$ws.Range("I4").Value = 0.09718402276460011
$ws.Range("J4").Value = 0.1059076069828809
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 9.964499267882889
$ws.Range("R4").Value = 89.680493410946
$ws.Range("S4").Value = 0.0001774293547557516
$ws.Range("T4").Value = 0.0002129292978800581
$ws.Range("G5").Value = 63.73255033333334
$ws.Range("H5").Value = 191.197651
$ws.Range("I5").Value = 0.09718402276460011
$ws.Range("J5").Value = 0.1059076069828809
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 1505.128908547335
$ws.Range("R5").Value = 9030.773451284009
$ws.Range("S5").Value = 0.02680054901790585
$ws.Range("T5").Value = 0.02144185627396504
$ws.Range("G6").Value = 63.73255033333334
$ws.Range("H6").Value = 191.197651
$ws.Range("I6").Value = 0.09718402276460011
$ws.Range("J6").Value = 0.1059076069828809
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 9.226561311756667
$ws.Range("R6").Value = 83.03905180581
$ws.Range("S6").Value = 0.0001642895218464084
$ws.Range("T6").Value = 0.0001971604562501077
$ws.Range("I7").Value = 0.1912449004891238
$ws.Range("J7").Value = 0.2084117242969288
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 82.75182797873077
$ws.Range("R7").Value = 744.766451808577
$ws.Range("S7").Value = 0.001473491346469304
$ws.Range("T7").Value = 0.001768306480446578
$ws.Range("I8").Value = 0.1912449004891238
$ws.Range("J8").Value = 0.2084117242969288
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("S8").Value = 0.1363591282463639
$ws.Range("T8").Value = 0.1636417687310214
$ws.Range("I9").Value = 0.1912449004891238
$ws.Range("J9").Value = 0.2084117242969288
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 19.60877535936244
$ws.Range("R9").Value = 176.478978234262
$ws.Range("S9").Value = 0.0003491567680451407
$ws.Range("T9").Value = 0.0004190158137714362
$ws.Range("I10").Value = 0.1912449004891238
$ws.Range("J10").Value = 0.2084117242969288
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 2961.888386074187
$ws.Range("R10").Value = 17771.33031644512
$ws.Range("S10").Value = 0.05273982475903715
$ws.Range("T10").Value = 0.04219464838730903
$ws.Range("I11").Value = 0.1912449004891238
$ws.Range("J11").Value = 0.2084117242969288
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 18.15661411956334
$ws.Range("R11").Value = 163.40952707607
$ws.Range("S11").Value = 0.0003232993692083187
$ws.Range("T11").Value = 0.0003879848843803645
$ws.Range("G12").Value = 184.1540323333334
$ws.Range("H12").Value = 552.4620970000001
$ws.Range("I12").Value = 0.2808114468489298
$ws.Range("J12").Value = 0.3060180830465028
$ws.Range("M12").Value = 0.6598136666666666
$ws.Range("N12").Value = 1.979441
$ws.Range("O12").Value = 0.007704735356083927
$ws.Range("P12").Value = 0.008484678519943686
$ws.Range("Q12").Value = 121.5073473053086
$ws.Range("R12").Value = 1093.566125747777
$ws.Range("S12").Value = 0.002163577882930032
$ws.Range("T12").Value = 0.002596465055939005
$ws.Range("G13").Value = 184.1540323333334
$ws.Range("H13").Value = 552.4620970000001
$ws.Range("I13").Value = 0.2808114468489298
$ws.Range("J13").Value = 0.3060180830465028
$ws.Range("O13").Value = 0.7130079175842846
$ws.Range("P13").Value = 0.7851850431306702
$ws.Range("Q13").Value = 11244.47455614915
$ws.Range("R13").Value = 101200.2710053424
$ws.Range("S13").Value = 0.2002207849515854
$ws.Range("T13").Value = 0.2402808217356333
$ws.Range("G14").Value = 184.1540323333334
$ws.Range("H14").Value = 552.4620970000001
$ws.Range("I14").Value = 0.2808114468489298
$ws.Range("J14").Value = 0.3060180830465028
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1563486666666667
$ws.Range("N14").Value = 0.469046
$ws.Range("O14").Value = 0.001825704984300993
$ws.Range("P14").Value = 0.002010519394650058
$ws.Range("Q14").Value = 28.79223741660689
$ws.Range("R14").Value = 259.130136749462
$ws.Range("S14").Value = 0.0005126788581608644
$ws.Range("T14").Value = 0.0006152552910786261
$ws.Range("G15").Value = 184.1540323333334
$ws.Range("H15").Value = 552.4620970000001
$ws.Range("I15").Value = 0.2808114468489298
$ws.Range("J15").Value = 0.3060180830465028
$ws.Range("M15").Value = 23.6163295
$ws.Range("N15").Value = 47.232659
$ws.Range("O15").Value = 0.2757711427815902
$ws.Range("P15").Value = 0.2024581319964196
$ws.Range("Q15").Value = 4349.042306337655
$ws.Range("R15").Value = 26094.25383802592
$ws.Range("S15").Value = 0.07743969360368112
$ws.Range("T15").Value = 0.06195584945072016
$ws.Range("G16").Value = 184.1540323333334
$ws.Range("H16").Value = 552.4620970000001
$ws.Range("I16").Value = 0.2808114468489298
$ws.Range("J16").Value = 0.3060180830465028
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.14477
$ws.Range("N16").Value = 0.43431
$ws.Range("O16").Value = 0.00169049929374041
$ws.Range("P16").Value = 0.001861626958316384
$ws.Range("Q16").Value = 26.65997926089667
$ws.Range("R16").Value = 239.9398133480701
$ws.Range("S16").Value = 0.0004747115525723384
$ws.Range("T16").Value = 0.0005696915131316718
$ws.Range("G17").Value = 162.052406
$ws.Range("H17").Value = 324.104812
$ws.Range("I17").Value = 0.2471092813859239
$ws.Range("J17").Value = 0.1795271274047008
$ws.Range("M17").Value = 0.6598136666666666
$ws.Range("N17").Value = 1.979441
$ws.Range("O17").Value = 0.007704735356083927
$ws.Range("P17").Value = 0.008484678519943686
$ws.Range("Q17").Value = 106.9243921950153
$ws.Range("R17").Value = 641.5463531700921
$ws.Range("S17").Value = 0.001903911617110619
$ws.Range("T17").Value = 0.001523229961637858
$ws.Range("G18").Value = 162.052406
$ws.Range("H18").Value = 324.104812
$ws.Range("I18").Value = 0.2471092813859239
$ws.Range("J18").Value = 0.1795271274047008
$ws.Range("O18").Value = 0.7130079175842846
$ws.Range("P18").Value = 0.7851850431306702
$ws.Range("Q18").Value = 9894.945730710024
$ws.Range("R18").Value = 59369.67438426014
$ws.Range("S18").Value = 0.1761908741367266
$ws.Range("T18").Value = 0.1409620152743853
$ws.Range("G19").Value = 162.052406
$ws.Range("H19").Value = 324.104812
$ws.Range("I19").Value = 0.2471092813859239
$ws.Range("J19").Value = 0.1795271274047008
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.1563486666666667
$ws.Range("N19").Value = 0.469046
$ws.Range("O19").Value = 0.001825704984300993
$ws.Range("P19").Value = 0.002010519394650058
$ws.Range("Q19").Value = 25.33667760822534
$ws.Range("R19").Value = 152.020065649352
$ws.Range("S19").Value = 0.0004511486466933178
$ws.Range("T19").Value = 0.000360942771512963
$ws.Range("G20").Value = 162.052406
$ws.Range("H20").Value = 324.104812
$ws.Range("I20").Value = 0.2471092813859239
$ws.Range("J20").Value = 0.1795271274047008
$ws.Range("M20").Value = 23.6163295
$ws.Range("N20").Value = 47.232659
$ws.Range("O20").Value = 0.2757711427815902
$ws.Range("P20").Value = 0.2024581319964196
$ws.Range("Q20").Value = 3827.083016363777
$ws.Range("R20").Value = 15308.33206545511
$ws.Range("S20").Value = 0.06814560891973374
$ws.Range("T20").Value = 0.03634672685703896
$ws.Range("G21").Value = 162.052406
$ws.Range("H21").Value = 324.104812
$ws.Range("I21").Value = 0.2471092813859239
$ws.Range("J21").Value = 0.1795271274047008
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.14477
$ws.Range("N21").Value = 0.43431
$ws.Range("O21").Value = 0.00169049929374041
$ws.Range("P21").Value = 0.001861626958316384
$ws.Range("Q21").Value = 23.46032681662
$ws.Range("R21").Value = 140.76196089972
$ws.Range("S21").Value = 0.0004177380656596045
$ws.Range("T21").Value = 0.0003342125401256912
$ws.Range("G22").Value = 120.4365156666667
$ws.Range("H22").Value = 361.309547
$ws.Range("I22").Value = 0.1836503485114226
$ws.Range("J22").Value = 0.2001354582689865
$ws.Range("M22").Value = 0.6598136666666666
$ws.Range("N22").Value = 1.979441
$ws.Range("O22").Value = 0.007704735356083927
$ws.Range("P22").Value = 0.008484678519943686
$ws.Range("Q22").Value = 79.46565900258078
$ws.Range("R22").Value = 715.190931023227
$ws.Range("S22").Value = 0.001414977333333093
$ws.Range("T22").Value = 0.001698085023853956
$ws.Range("G23").Value = 120.4365156666667
$ws.Range("H23").Value = 361.309547
$ws.Range("I23").Value = 0.1836503485114226
$ws.Range("J23").Value = 0.2001354582689865
$ws.Range("O23").Value = 0.7130079175842846
$ws.Range("P23").Value = 0.7851850431306702
$ws.Range("Q23").Value = 7353.872836158162
$ws.Range("R23").Value = 66184.85552542347
$ws.Range("S23").Value = 0.1309441525557576
$ws.Range("T23").Value = 0.1571433684329106
$ws.Range("G24").Value = 120.4365156666667
$ws.Range("H24").Value = 361.309547
$ws.Range("I24").Value = 0.1836503485114226
$ws.Range("J24").Value = 0.2001354582689865
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1563486666666667
$ws.Range("N24").Value = 0.469046
$ws.Range("O24").Value = 0.001825704984300993
$ws.Range("P24").Value = 0.002010519394650058
$ws.Range("Q24").Value = 18.83008864246244
$ws.Range("R24").Value = 169.470797782162
$ws.Range("S24").Value = 0.0003352913566459187
$ws.Range("T24").Value = 0.0004023762204069749
$ws.Range("G25").Value = 120.4365156666667
$ws.Range("H25").Value = 361.309547
$ws.Range("I25").Value = 0.1836503485114226
$ws.Range("J25").Value = 0.2001354582689865
$ws.Range("M25").Value = 23.6163295
$ws.Range("N25").Value = 47.232659
$ws.Range("O25").Value = 0.2757711427815902
$ws.Range("P25").Value = 0.2024581319964196
$ws.Range("Q25").Value = 2844.268437815912
$ws.Range("R25").Value = 17065.61062689547
$ws.Range("S25").Value = 0.05064546648123232
$ws.Range("T25").Value = 0.04051905102738641
$ws.Range("G26").Value = 120.4365156666667
$ws.Range("H26").Value = 361.309547
$ws.Range("I26").Value = 0.1836503485114226
$ws.Range("J26").Value = 0.2001354582689865
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.14477
$ws.Range("N26").Value = 0.43431
$ws.Range("O26").Value = 0.00169049929374041
$ws.Range("P26").Value = 0.001861626958316384
$ws.Range("Q26").Value = 17.43559437306333
$ws.Range("R26").Value = 156.92034935757
$ws.Range("S26").Value = 0.0003104607844537401
$ws.Range("T26").Value = 0.0003725775644285491
